$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.7
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 4.2
$ws.Range("J2").Value = 2.2
$ws.Range("L2").Value = 4.33
$ws.Range("N2").Value = 21
$ws.Range("U2").Value = 1.44
$ws.Range("V2").Value = 2.63
$ws.Range("AD2").Value = 8.5
$ws.Range("AE2").Value = 12
$ws.Range("AH2").Value = 19
$ws.Range("AI2").Value = 26
$ws.Range("AK2").Value = 41
$ws.Range("AM2").Value = 26
$ws.Range("AN2").Value = 4.33
$ws.Range("AO2").Value = 8.5
$ws.Range("AQ2").Value = 23
$ws.Range("AR2").Value = 41
$ws.Range("AX2").Value = 6.5
$ws.Range("AZ2").Value = 21
$ws.Range("BA2").Value = 51

# Row 4 updates
$ws.Range("G4").Value = 2.5
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 2.9
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 3.75
$ws.Range("N4").Value = 7.5
$ws.Range("X4").Value = 11
$ws.Range("Y4").Value = 10
$ws.Range("Z4").Value = 23
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 7.5
$ws.Range("AN4").Value = 4.33

# Row 6 updates
$ws.Range("G6").Value = 1.9
$ws.Range("I6").Value = 4
$ws.Range("Q6").Value = 1.88
$ws.Range("R6").Value = 1.98
$ws.Range("Z6").Value = 17
$ws.Range("AD6").Value = 6.5
$ws.Range("AL6").Value = 29
$ws.Range("AU6").Value = 7.5
$ws.Range("AZ6").Value = 26
